$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '30.744.90'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +2.05%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.115.28'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +10.28%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '333.40'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +4.15%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5236'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +3.54%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.4424'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +8.51%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.09059'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +8.45%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '46.29'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +9.17%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '1.181'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +6.44%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '25.24'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +4.86%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '2.114.45'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +11.21%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '6.826'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +6.43%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '7.753'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +6.99%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '98.50'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +6.44%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.32%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.00001138'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +4.01%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.06658'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +2.29%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '19.28'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +3.98%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.411'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +7.84%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '30.859.71'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +2.38%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '12.06'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +6.22%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.363.90'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +11.04%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.259'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +3.03%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '22.96'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +5.13%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.552'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +11.81%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '163.39'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.31%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '133.92'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +4.01%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.183'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +3.45%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.1070'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +2.41%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '6.255'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +5.12%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.540'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +28.42%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '3.912'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +3.50%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.02589'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +5.36%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '5.619'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +4.71%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.06781'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +5.18%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '9.583'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +11.08%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '12.79'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +12.08%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.2274'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +5.43%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.6796'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +4.05%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.255'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +3.56%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '14.25'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +6.47%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.9996'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.02%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.6342'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +4.30%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '2.254'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +2.92%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.294'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +6.88%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '3.670'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +1.28%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '83.24'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +5.25%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '121.56'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.51%  '
